$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("sim") values recalibrated after adding work dummy to utility fn
$ws.Range("C3").Value = 0.653232397092348
$ws.Range("C4").Value = 0.494567508279926
$ws.Range("C5").Value = -0.1970354481462937
$ws.Range("C6").Value = -0.22964324383684
$ws.Range("C7").Value = 0.08201177580355344
$ws.Range("C8").Value = 0.1776880348197343
$ws.Range("C9").Value = 0.008220003613187324
$ws.Range("C10").Value = 0.0001622621681291687
$ws.Range("C11").Value = -0.005741201080244457
$ws.Range("C12").Value = 1.440248505367899

# Minor downstream floating point shifts in "data"/"SE" columns from recalibration
$ws.Range("D6").Value = -0.3254977791559274
$ws.Range("E7").Value = 0.002933956643042229
$ws.Range("E8").Value = 0.00006191117344771939
$ws.Range("D9").Value = 0.1611544252757249
$ws.Range("E9").Value = 0.02682390016553195
$ws.Range("D10").Value = -0.003200531431248727
